# cryptos.xlsx refresh -- GitHub Actions scheduled price/volume update.
# Re-scrapes coinranking.com and rewrites the Coin / Link / Price / Volume(1h)
# columns (B:E) on the active sheet. Row 48-50 also picked up a ranking reshuffle
# (Quant / NEARProtocol / PEPE rotated), so Coin + Link are rewritten there too.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "28.972.21"
$ws.Range("E2").Value = "  -1.58%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "1.906.29"
$ws.Range("E3").Value = "  -3.10%  "

# Row 4: TetherUSD
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.22%  "

# Row 5: BNB
$ws.Range("D5").Value = "'324.29"
$ws.Range("E5").Value = "  -0.59%  "

# Row 6: USDC
$ws.Range("E6").Value = "  +0.02%  "

# Row 7: XRP
$ws.Range("D7").Value = "'0.4594"
$ws.Range("E7").Value = "  -1.26%  "

# Row 8: Cardano
$ws.Range("D8").Value = "'0.3828"
$ws.Range("E8").Value = "  -2.08%  "

# Row 9: Dogecoin
$ws.Range("D9").Value = "'0.07715"
$ws.Range("E9").Value = "  -2.59%  "

# Row 10: Polygon
$ws.Range("D10").Value = "'0.9797"
$ws.Range("E10").Value = "  -0.74%  "

# Row 11: Solana
$ws.Range("D11").Value = "'22.06"
$ws.Range("E11").Value = "  -2.97%  "

# Row 12: WrappedEther
$ws.Range("D12").Value = "1.901.64"
$ws.Range("E12").Value = "  -3.37%  "

# Row 13: Chainlink
$ws.Range("D13").Value = "'6.936"
$ws.Range("E13").Value = "  -3.25%  "

# Row 14: Polkadot
$ws.Range("D14").Value = "'5.659"
$ws.Range("E14").Value = "  -2.63%  "

# Row 15: TRON
$ws.Range("D15").Value = "'0.07044"
$ws.Range("E15").Value = "  -0.53%  "

# Row 16: BinanceUSD
$ws.Range("E16").Value = "  +0.09%  "

# Row 17: Litecoin
$ws.Range("D17").Value = "'83.77"
$ws.Range("E17").Value = "  -4.46%  "

# Row 18: ShibaInu
$ws.Range("D18").Value = "'0.000009480"
$ws.Range("E18").Value = "  -4.34%  "

# Row 19: Avalanche
$ws.Range("D19").Value = "'16.65"
$ws.Range("E19").Value = "  -3.43%  "

# Row 20: Dai
$ws.Range("E20").Value = "  -0.03%  "

# Row 21: WrappedBTC
$ws.Range("D21").Value = "28.988.78"
$ws.Range("E21").Value = "  -1.47%  "

# Row 22: Uniswap
$ws.Range("D22").Value = "'5.303"

# Row 23: Cosmos
$ws.Range("D23").Value = "'10.88"
$ws.Range("E23").Value = "  -2.09%  "

# Row 24: Toncoin
$ws.Range("D24").Value = "'2.092"
$ws.Range("E24").Value = "  -0.60%  "

# Row 25: Monero
$ws.Range("D25").Value = "'158.00"
$ws.Range("E25").Value = "  +0.00%  "

# Row 26: EthereumClassic
$ws.Range("D26").Value = "'19.06"
$ws.Range("E26").Value = "  -2.06%  "

# Row 27: InternetComputer(DFINITY)
$ws.Range("D27").Value = "'5.648"
$ws.Range("E27").Value = "  -2.16%  "

# Row 28: BitcoinCash
$ws.Range("D28").Value = "'117.45"
$ws.Range("E28").Value = "  -1.57%  "

# Row 29: LidoDAOToken
$ws.Range("D29").Value = "'1.850"
$ws.Range("E29").Value = "  -2.43%  "

# Row 30: Stellar
$ws.Range("D30").Value = "'0.09279"
$ws.Range("E30").Value = "  -1.35%  "

# Row 31: ImmutableX
$ws.Range("D31").Value = "'0.8665"
$ws.Range("E31").Value = "  -2.61%  "

# Row 32: Filecoin
$ws.Range("D32").Value = "'5.080"
$ws.Range("E32").Value = "  -2.82%  "

# Row 33: ARBITRUM
$ws.Range("D33").Value = "'1.248"
$ws.Range("E33").Value = "  -5.18%  "

# Row 34: HuobiToken
$ws.Range("D34").Value = "'2.944"
$ws.Range("E34").Value = "  -6.99%  "

# Row 35: Hedera
$ws.Range("D35").Value = "'0.05719"
$ws.Range("E35").Value = "  -1.46%  "

# Row 36: TrustWalletToken
$ws.Range("D36").Value = "'1.145"
$ws.Range("E36").Value = "  -1.76%  "

# Row 37: Frax
$ws.Range("E37").Value = "  +0.09%  "

# Row 38: VeChain
$ws.Range("D38").Value = "'0.02041"
$ws.Range("E38").Value = "  -2.69%  "

# Row 39: TheSandbox
$ws.Range("D39").Value = "'0.5502"
$ws.Range("E39").Value = "  -3.51%  "

# Row 40: FraxShare
$ws.Range("D40").Value = "'7.405"
$ws.Range("E40").Value = "  -4.16%  "

# Row 41: Algorand
$ws.Range("D41").Value = "'0.1756"
$ws.Range("E41").Value = "  -2.04%  "

# Row 42: Aptos
$ws.Range("D42").Value = "'9.316"
$ws.Range("E42").Value = "  -3.34%  "

# Row 43: MXToken
$ws.Range("D43").Value = "'2.776"
$ws.Range("E43").Value = "  +1.04%  "

# Row 44: Decentraland
$ws.Range("D44").Value = "'0.5178"
$ws.Range("E44").Value = "  -2.87%  "

# Row 45: EnergySwap
$ws.Range("D45").Value = "'11.25"
$ws.Range("E45").Value = "  -4.38%  "

# Row 46: Cronos
$ws.Range("D46").Value = "'0.06836"
$ws.Range("E46").Value = "  -1.10%  "

# Row 47: RenderToken
$ws.Range("E47").Value = "  -5.21%  "

# Row 48: Quant -> PEPE
$ws.Range("B48").Value = "PEPE"
$ws.Range("C48").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D48").Value = "'0.000002582"
$ws.Range("E48").Value = "  -17.57%  "

# Row 49: NEARProtocol -> Quant
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "'110.50"
$ws.Range("E49").Value = "  -2.38%  "

# Row 50: PEPE -> NEARProtocol
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'1.776"
$ws.Range("E50").Value = "  -2.63%  "

# Row 51: PaxDollar
$ws.Range("E51").Value = "  -0.04%  "
